# Update all data rows (2-17) with refreshed TPM-based NATMI metrics,
# and append the new "Resolving-Mac" sending-cluster rows (14-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
#          E..T  the 16 numeric NATMI metrics (same order as header row 1)
$rows = @(
    @("ECs", "Mmrn2", "Cd248", "ECs", 3.0, 1.0, 46.105436, 138.316308, 0.929514990096053, 0.9295149900960532, 3.0, 1.0, 1.765932333333333, 5.297797, 0.02646494741755958, 0.02646494741755958, 81.41908017483067, 732.7717215734759, 0.02459956533672546, 0.02459956533672546),
    @("ECs", "Mmrn2", "Cd248", "FAPs", 3.0, 1.0, 46.105436, 138.316308, 0.929514990096053, 0.9295149900960532, 3.0, 1.0, 59.98702866666667, 179.961086, 0.8989888963652094, 0.8989888963652094, 2765.728111021166, 24891.55299919049, 0.8356236551013693, 0.8356236551013694),
    @("ECs", "Mmrn2", "Cd248", "MuSCs", 3.0, 1.0, 46.105436, 138.316308, 0.929514990096053, 0.9295149900960532, 3.0, 1.0, 4.969642, 14.908926, 0.07447698404493167, 0.07447698404493167, 229.127511173912, 2062.147600565208, 0.06922747308690855, 0.06922747308690858),
    @("ECs", "Mmrn2", "Cd248", "Resolving-Mac", 3.0, 1.0, 46.105436, 138.316308, 0.929514990096053, 0.9295149900960532, 1.0, 0.3333333333333333, 0.004615666666666666, 0.013847, 0.00006917217229934394, 0.00006917217229934394, 0.2128073240973333, 1.915265916876, 0.00006429657104974714, 0.00006429657104974716),
    @("FAPs", "Mmrn2", "Cd248", "ECs", 3.0, 1.0, 2.068862333333333, 6.206586999999999, 0.04170958390412858, 0.04170958390412859, 3.0, 1.0, 1.765932333333333, 5.297797, 0.02646494741755958, 0.02646494741755958, 3.653470887648778, 32.881237988839, 0.001103841944831052, 0.001103841944831052),
    @("FAPs", "Mmrn2", "Cd248", "FAPs", 3.0, 1.0, 2.068862333333333, 6.206586999999999, 0.04170958390412858, 0.04170958390412859, 3.0, 1.0, 59.98702866666667, 179.961086, 0.8989888963652094, 0.8989888963652094, 124.1049040970536, 1116.944136873482, 0.03749645280182465, 0.03749645280182466),
    @("FAPs", "Mmrn2", "Cd248", "MuSCs", 3.0, 1.0, 2.068862333333333, 6.206586999999999, 0.04170958390412858, 0.04170958390412859, 3.0, 1.0, 4.969642, 14.908926, 0.07447698404493167, 0.07447698404493167, 10.28150514395133, 92.53354629556199, 0.003106404014948523, 0.003106404014948523),
    @("FAPs", "Mmrn2", "Cd248", "Resolving-Mac", 3.0, 1.0, 2.068862333333333, 6.206586999999999, 0.04170958390412858, 0.04170958390412859, 1.0, 0.3333333333333333, 0.004615666666666666, 0.013847, 0.00006917217229934394, 0.00006917217229934394, 0.009549178909888887, 0.08594261018899998, 0.000002885142524350325, 0.000002885142524350325),
    @("MuSCs", "Mmrn2", "Cd248", "ECs", 3.0, 1.0, 1.378131333333333, 4.134394, 0.02778400648145685, 0.02778400648145685, 3.0, 1.0, 1.765932333333333, 5.297797, 0.02646494741755958, 0.02646494741755958, 2.433686681113111, 21.903180130018, 0.00073530227058089, 0.0007353022705808901),
    @("MuSCs", "Mmrn2", "Cd248", "FAPs", 3.0, 1.0, 1.378131333333333, 4.134394, 0.02778400648145685, 0.02778400648145685, 3.0, 1.0, 59.98702866666667, 179.961086, 0.8989888963652094, 0.8989888963652094, 82.67000379909824, 744.0300341918842, 0.02497751332336872, 0.02497751332336872),
    @("MuSCs", "Mmrn2", "Cd248", "MuSCs", 3.0, 1.0, 1.378131333333333, 4.134394, 0.02778400648145685, 0.02778400648145685, 3.0, 1.0, 4.969642, 14.908926, 0.07447698404493167, 0.07447698404493167, 6.848819355649334, 61.63937420084401, 0.002069269007423739, 0.00206926900742374),
    @("MuSCs", "Mmrn2", "Cd248", "Resolving-Mac", 3.0, 1.0, 1.378131333333333, 4.134394, 0.02778400648145685, 0.02778400648145685, 1.0, 0.3333333333333333, 0.004615666666666666, 0.013847, 0.00006917217229934394, 0.00006917217229934394, 0.006360994857555555, 0.057248953718, 0.000001921880083501422, 0.000001921880083501422),
    @("Resolving-Mac", "Mmrn2", "Cd248", "ECs", 1.0, 0.3333333333333333, 0.049176, 0.147528, 0.0009914195183614249, 0.000991419518361425, 3.0, 1.0, 1.765932333333333, 5.297797, 0.02646494741755958, 0.02646494741755958, 0.086841488424, 0.7815733958159999, 0.00002623786542217735, 0.00002623786542217735),
    @("Resolving-Mac", "Mmrn2", "Cd248", "FAPs", 1.0, 0.3333333333333333, 0.049176, 0.147528, 0.0009914195183614249, 0.000991419518361425, 3.0, 1.0, 59.98702866666667, 179.961086, 0.8989888963652094, 0.8989888963652094, 2.949922121712, 26.549299095408, 0.0008912751386466648, 0.000891275138646665),
    @("Resolving-Mac", "Mmrn2", "Cd248", "MuSCs", 1.0, 0.3333333333333333, 0.049176, 0.147528, 0.0009914195183614249, 0.000991419518361425, 3.0, 1.0, 4.969642, 14.908926, 0.07447698404493167, 0.07447698404493167, 0.244387114992, 2.199484034928, 0.00007383793565083767, 0.0000738379356508377),
    @("Resolving-Mac", "Mmrn2", "Cd248", "Resolving-Mac", 1.0, 0.3333333333333333, 0.049176, 0.147528, 0.0009914195183614249, 0.000991419518361425, 1.0, 0.3333333333333333, 0.004615666666666666, 0.013847, 0.00006917217229934394, 0.00006917217229934394, 0.000226980024, 0.002042820216, 0.00000006857864174502906, 0.00000006857864174502907)
)

$numRows = $rows.Count
$numCols = 20
$data = New-Object 'object[,]' $numRows,$numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("A2:T17").Value = $data

